$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ASN")

# Update QTY (column H) values for the specified lines now that
# the Chewy label and ASN work is complete.
# The leading apostrophe forces Excel to keep these numeric-looking
# entries as text, matching the workbook's existing text-typed QTY cells.
# Re-applying the "Normal" style afterward keeps the cell formatting
# identical to what it was before the edit.
$ws.Range("H21").Value = "'4"
$ws.Range("H21").Style = "Normal"

$ws.Range("H23").Value = "'5"
$ws.Range("H23").Style = "Normal"

$ws.Range("H26").Value = "'2"
$ws.Range("H26").Style = "Normal"

$ws.Range("H28").Value = "'6"
$ws.Range("H28").Style = "Normal"
